# fhir ig initial setup
# Regenerated FHIR IG export: bump the generation timestamp and drop the
# three SNOMED concept rows (55929007/418107008/214264003 - "Feeling
# irritable (finding)" / "Unconscious (finding)" / "Lethargy (finding)")
# that are no longer part of the included value set.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the generation timestamp -------------------------
$metaWs = $wb.Worksheets.Item("Metadata")
$metaWs.Range("B8").Value = "2025-08-20T08:30:34+05:45"

# --- Include #0 sheet: remove rows 5-7 (55929007/418107008/214264003) ------
# and shift the remaining rows (the blank separator row + "System URI" row)
# up so they land on rows 5-6.
$incWs = $wb.Worksheets.Item("Include #0")
$incWs.Rows("5:7").Delete()
